$wb = $excel.ActiveWorkbook

$ws_LP1912 = $wb.Worksheets.Item("LP1912")
$ws_LP1912.Cells.Item(2, 1).Value = "Última actualización: 16:45:31"
$ws_LP1912.Cells.Item(3, 1).Value = "Total filas: 354"
$ws_LP1912.Cells.Item(15, 1).Value = "03:45:25"
$ws_LP1912.Cells.Item(15, 3).Value = "215B_EL PATO"
$ws_LP1912.Cells.Item(15, 4).Value = 109
$ws_LP1912.Cells.Item(16, 1).Value = "04:18:02"
$ws_LP1912.Cells.Item(16, 3).Value = "14_ABASTO"
$ws_LP1912.Cells.Item(16, 4).Value = 76
$ws_LP1912.Cells.Item(70, 1).Value = "07:19:29"
$ws_LP1912.Cells.Item(70, 3).Value = "11_ETCHEVERRY"
$ws_LP1912.Cells.Item(70, 4).Value = 44
$ws_LP1912.Cells.Item(71, 1).Value = "06:55:02"
$ws_LP1912.Cells.Item(71, 3).Value = "23_HERNANDEZ"
$ws_LP1912.Cells.Item(71, 4).Value = 68
$ws_LP1912.Cells.Item(79, 1).Value = "06:55:02"
$ws_LP1912.Cells.Item(79, 3).Value = "16_P MOR-SANTA ANA"
$ws_LP1912.Cells.Item(79, 4).Value = 88
$ws_LP1912.Cells.Item(80, 1).Value = "06:25:43"
$ws_LP1912.Cells.Item(80, 3).Value = "215B_EL PATO"
$ws_LP1912.Cells.Item(80, 4).Value = 118
$ws_LP1912.Cells.Item(103, 1).Value = "08:32:09"
$ws_LP1912.Cells.Item(103, 3).Value = "16_SANTA ANA"
$ws_LP1912.Cells.Item(103, 4).Value = 51
$ws_LP1912.Cells.Item(105, 1).Value = "08:02:22"
$ws_LP1912.Cells.Item(105, 3).Value = "11_ETCHEVERRY"
$ws_LP1912.Cells.Item(105, 4).Value = 81
$ws_LP1912.Cells.Item(204, 1).Value = "10:59:49"
$ws_LP1912.Cells.Item(204, 3).Value = "11_ETCHEVERRY"
$ws_LP1912.Cells.Item(204, 4).Value = 109
$ws_LP1912.Cells.Item(205, 1).Value = "12:47:27"
$ws_LP1912.Cells.Item(205, 3).Value = "16_SANTA ANA"
$ws_LP1912.Cells.Item(205, 4).Value = 1
$ws_LP1912.Cells.Item(215, 1).Value = "12:21:08"
$ws_LP1912.Cells.Item(215, 3).Value = "16_P MOR-SANTA ANA"
$ws_LP1912.Cells.Item(215, 4).Value = 46
$ws_LP1912.Cells.Item(216, 1).Value = "11:30:45"
$ws_LP1912.Cells.Item(216, 3).Value = "10_OLMOS"
$ws_LP1912.Cells.Item(216, 4).Value = 97
$ws_LP1912.Cells.Item(224, 3).Value = "15_ABASTO"
$ws_LP1912.Cells.Item(225, 3).Value = "14_ABASTO"
$ws_LP1912.Cells.Item(226, 1).Value = "11:56:55"
$ws_LP1912.Cells.Item(226, 3).Value = "10_OLMOS"
$ws_LP1912.Cells.Item(226, 4).Value = 91
$ws_LP1912.Cells.Item(227, 1).Value = "12:21:08"
$ws_LP1912.Cells.Item(227, 3).Value = "14_ABASTO"
$ws_LP1912.Cells.Item(227, 4).Value = 66
$ws_LP1912.Cells.Item(248, 1).Value = "12:21:08"
$ws_LP1912.Cells.Item(248, 3).Value = "17_ROMERO"
$ws_LP1912.Cells.Item(248, 4).Value = 103
$ws_LP1912.Cells.Item(249, 1).Value = "13:33:42"
$ws_LP1912.Cells.Item(249, 3).Value = "23_HERNANDEZ"
$ws_LP1912.Cells.Item(249, 4).Value = 31
$ws_LP1912.Cells.Item(292, 1).Value = "14:24:16"
$ws_LP1912.Cells.Item(292, 3).Value = "14_ABASTO"
$ws_LP1912.Cells.Item(292, 4).Value = 82
$ws_LP1912.Cells.Item(293, 1).Value = "14:56:20"
$ws_LP1912.Cells.Item(293, 3).Value = "16_P MOR-167 Y 521"
$ws_LP1912.Cells.Item(293, 4).Value = 50
$ws_LP1912.Cells.Item(327, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(327, 2).Value = "16:45"
$ws_LP1912.Cells.Item(327, 3).Value = "16_SANTA ANA"
$ws_LP1912.Cells.Item(327, 4).Value = 0
$ws_LP1912.Cells.Item(328, 1).Value = "15:22:17"
$ws_LP1912.Cells.Item(328, 2).Value = "16:48"
$ws_LP1912.Cells.Item(328, 3).Value = "15_ABASTO"
$ws_LP1912.Cells.Item(328, 4).Value = 86
$ws_LP1912.Cells.Item(329, 1).Value = "15:53:28"
$ws_LP1912.Cells.Item(329, 2).Value = "16:50"
$ws_LP1912.Cells.Item(329, 3).Value = "14_ABASTO"
$ws_LP1912.Cells.Item(329, 4).Value = 57
$ws_LP1912.Cells.Item(330, 1).Value = "15:22:17"
$ws_LP1912.Cells.Item(330, 3).Value = "17_179 Y 38"
$ws_LP1912.Cells.Item(330, 4).Value = 94
$ws_LP1912.Cells.Item(331, 1).Value = "16:31:51"
$ws_LP1912.Cells.Item(331, 2).Value = "16:56"
$ws_LP1912.Cells.Item(331, 4).Value = 25
$ws_LP1912.Cells.Item(332, 1).Value = "16:13:37"
$ws_LP1912.Cells.Item(332, 2).Value = "16:57"
$ws_LP1912.Cells.Item(332, 3).Value = "10_OLMOS"
$ws_LP1912.Cells.Item(332, 4).Value = 44
$ws_LP1912.Cells.Item(333, 1).Value = "15:22:17"
$ws_LP1912.Cells.Item(333, 3).Value = "215A_EL PATO"
$ws_LP1912.Cells.Item(333, 4).Value = 102
$ws_LP1912.Cells.Item(334, 3).Value = "11_ETCHEVERRY"
$ws_LP1912.Cells.Item(335, 1).Value = "16:13:37"
$ws_LP1912.Cells.Item(335, 2).Value = "17:04"
$ws_LP1912.Cells.Item(335, 3).Value = "23_HERNANDEZ"
$ws_LP1912.Cells.Item(335, 4).Value = 51
$ws_LP1912.Cells.Item(336, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(336, 2).Value = "17:06"
$ws_LP1912.Cells.Item(336, 3).Value = "23_HERNANDEZ"
$ws_LP1912.Cells.Item(336, 4).Value = 21
$ws_LP1912.Cells.Item(337, 2).Value = "17:09"
$ws_LP1912.Cells.Item(337, 3).Value = "10_OLMOS"
$ws_LP1912.Cells.Item(337, 4).Value = 38
$ws_LP1912.Cells.Item(338, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(338, 2).Value = "17:10"
$ws_LP1912.Cells.Item(338, 3).Value = "10_OLMOS"
$ws_LP1912.Cells.Item(338, 4).Value = 25
$ws_LP1912.Cells.Item(339, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(339, 2).Value = "17:16"
$ws_LP1912.Cells.Item(339, 3).Value = "11_ETCHEVERRY"
$ws_LP1912.Cells.Item(339, 4).Value = 31
$ws_LP1912.Cells.Item(340, 1).Value = "16:31:51"
$ws_LP1912.Cells.Item(340, 2).Value = "17:20"
$ws_LP1912.Cells.Item(340, 3).Value = "26_HERNANDEZ"
$ws_LP1912.Cells.Item(340, 4).Value = 49
$ws_LP1912.Cells.Item(341, 1).Value = "16:31:51"
$ws_LP1912.Cells.Item(341, 2).Value = "17:20"
$ws_LP1912.Cells.Item(341, 3).Value = "16_SANTA ANA"
$ws_LP1912.Cells.Item(341, 4).Value = 49
$ws_LP1912.Cells.Item(342, 2).Value = "17:21"
$ws_LP1912.Cells.Item(342, 3).Value = "26_HERNANDEZ"
$ws_LP1912.Cells.Item(342, 4).Value = 88
$ws_LP1912.Cells.Item(343, 2).Value = "17:24"
$ws_LP1912.Cells.Item(343, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws_LP1912.Cells.Item(343, 4).Value = 91
$ws_LP1912.Cells.Item(344, 1).Value = "15:53:28"
$ws_LP1912.Cells.Item(344, 2).Value = "17:28"
$ws_LP1912.Cells.Item(344, 3).Value = "14_ABASTO"
$ws_LP1912.Cells.Item(344, 4).Value = 95
$ws_LP1912.Cells.Item(345, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(345, 2).Value = "17:34"
$ws_LP1912.Cells.Item(345, 3).Value = "23_HERNANDEZ"
$ws_LP1912.Cells.Item(345, 4).Value = 49
$ws_LP1912.Cells.Item(346, 2).Value = "17:36"
$ws_LP1912.Cells.Item(346, 3).Value = "27_EL RETIRO"
$ws_LP1912.Cells.Item(346, 4).Value = 103
$ws_LP1912.Cells.Item(347, 2).Value = "17:38"
$ws_LP1912.Cells.Item(347, 3).Value = "17_ROMERO"
$ws_LP1912.Cells.Item(347, 4).Value = 105
$ws_LP1912.Cells.Item(348, 1).Value = "15:53:28"
$ws_LP1912.Cells.Item(348, 2).Value = "17:40"
$ws_LP1912.Cells.Item(348, 3).Value = "215B_EL PATO"
$ws_LP1912.Cells.Item(348, 4).Value = 107
$ws_LP1912.Cells.Item(349, 1).Value = "16:13:37"
$ws_LP1912.Cells.Item(349, 2).Value = "17:40"
$ws_LP1912.Cells.Item(349, 3).Value = "17_ROMERO"
$ws_LP1912.Cells.Item(349, 4).Value = 87
$ws_LP1912.Cells.Item(350, 2).Value = "17:45"
$ws_LP1912.Cells.Item(350, 3).Value = "15_ABASTO"
$ws_LP1912.Cells.Item(350, 4).Value = 74
$ws_LP1912.Cells.Item(351, 1).Value = "15:53:28"
$ws_LP1912.Cells.Item(351, 2).Value = "17:50"
$ws_LP1912.Cells.Item(351, 3).Value = "16_P MOR-167 Y 521"
$ws_LP1912.Cells.Item(351, 4).Value = 117
$ws_LP1912.Cells.Item(351, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(352, 1).Value = "15:53:28"
$ws_LP1912.Cells.Item(352, 2).Value = "17:52"
$ws_LP1912.Cells.Item(352, 3).Value = "81_EL PELIGRO"
$ws_LP1912.Cells.Item(352, 4).Value = 119
$ws_LP1912.Cells.Item(352, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(353, 1).Value = "16:13:37"
$ws_LP1912.Cells.Item(353, 2).Value = "18:04"
$ws_LP1912.Cells.Item(353, 3).Value = "17_ROMERO"
$ws_LP1912.Cells.Item(353, 4).Value = 111
$ws_LP1912.Cells.Item(353, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(354, 1).Value = "16:31:51"
$ws_LP1912.Cells.Item(354, 2).Value = "18:20"
$ws_LP1912.Cells.Item(354, 3).Value = "26_HERNANDEZ"
$ws_LP1912.Cells.Item(354, 4).Value = 109
$ws_LP1912.Cells.Item(354, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(355, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(355, 2).Value = "18:21"
$ws_LP1912.Cells.Item(355, 3).Value = "26_HERNANDEZ"
$ws_LP1912.Cells.Item(355, 4).Value = 96
$ws_LP1912.Cells.Item(355, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(356, 1).Value = "16:31:51"
$ws_LP1912.Cells.Item(356, 2).Value = "18:27"
$ws_LP1912.Cells.Item(356, 3).Value = "215C_EL PATO"
$ws_LP1912.Cells.Item(356, 4).Value = 116
$ws_LP1912.Cells.Item(356, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(357, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(357, 2).Value = "18:28"
$ws_LP1912.Cells.Item(357, 3).Value = "215C_EL PATO"
$ws_LP1912.Cells.Item(357, 4).Value = 103
$ws_LP1912.Cells.Item(357, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(358, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(358, 2).Value = "18:32"
$ws_LP1912.Cells.Item(358, 3).Value = "11X44_ETCHEVERRY"
$ws_LP1912.Cells.Item(358, 4).Value = 107
$ws_LP1912.Cells.Item(358, 5).Value = "LP1912"
$ws_LP1912.Cells.Item(359, 1).Value = "16:45:31"
$ws_LP1912.Cells.Item(359, 2).Value = "18:40"
$ws_LP1912.Cells.Item(359, 3).Value = "15_ABASTO"
$ws_LP1912.Cells.Item(359, 4).Value = 115
$ws_LP1912.Cells.Item(359, 5).Value = "LP1912"

$ws_LP1912_215 = $wb.Worksheets.Item("LP1912-215")
$ws_LP1912_215.Cells.Item(2, 1).Value = "Última actualización: 16:45:31"
$ws_LP1912_215.Cells.Item(3, 1).Value = "Total filas: 38"
$ws_LP1912_215.Cells.Item(43, 1).Value = "16:45:31"
$ws_LP1912_215.Cells.Item(43, 2).Value = "18:28"
$ws_LP1912_215.Cells.Item(43, 3).Value = "215C_EL PATO"
$ws_LP1912_215.Cells.Item(43, 4).Value = 103
$ws_LP1912_215.Cells.Item(43, 5).Value = "LP1912"

$ws_6203_6173 = $wb.Worksheets.Item("6203-6173")
$ws_6203_6173.Cells.Item(2, 1).Value = "Última actualización: 16:45:31"
$ws_6203_6173.Cells.Item(3, 1).Value = "Total filas: 48"
$ws_6203_6173.Cells.Item(53, 1).Value = "16:45:31"
$ws_6203_6173.Cells.Item(53, 2).Value = "18:04"
$ws_6203_6173.Cells.Item(53, 3).Value = "215C_LA PLATA"
$ws_6203_6173.Cells.Item(53, 4).Value = 79
$ws_6203_6173.Cells.Item(53, 5).Value = "L6203"

